# CIERRE 24 NOV 2021
# Fills in the daily purchase/expense entries for 10-Nov-2021 .. 22-Nov-2021
# (rows 7-19) on the "NOVIEMBRE 2 0 2 1" sheet, and updates the selection
# to match where the author left off (Q19, scrolled to G7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 7  (10-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C7").Value = 80
$ws.Range("D7").Value = "QUESO"
$ws.Range("F7").Value = 28829
$ws.Range("I7").Value = 33
$ws.Range("M7").Formula = "=22800+5000"
$ws.Range("N7").Value = 917

# ---------------------------------------------------------------
# Row 8  (11-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C8").Value = 19923
$ws.Range("D8").Value = "TOCINO-JAMON-CONDIMENTOS-QUESO-PAN ARABE-RETAZO-"
$ws.Range("F8").Value = 49002
$ws.Range("I8").Value = 11
$ws.Range("M8").Value = 29070

# ---------------------------------------------------------------
# Row 9  (12-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C9").Value = 3480
$ws.Range("D9").Value = "CHORIZO"
$ws.Range("F9").Value = 78448
$ws.Range("I9").Value = 1659
$ws.Range("M9").Formula = "=25000+47260"
$ws.Range("N9").Value = 1050

# ---------------------------------------------------------------
# Row 10  (13-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C10").Value = 19
$ws.Range("D10").Value = "CEBOLLA--PEREJIL"
$ws.Range("F10").Value = 77465
$ws.Range("I10").Value = 129
$ws.Range("J10").Value = 44513
$ws.Range("K10").Value = "NOMINA # 46 y vac"
$ws.Range("L10").Formula = "=10814.29+1285.71"
$ws.Range("M10").Formula = "=30000+29340"
$ws.Range("N10").Value = 5875

# ---------------------------------------------------------------
# Row 11  (14-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C11").Value = 962
$ws.Range("D11").Value = "PAPAS-SANCHICHA-CECINA"
$ws.Range("F11").Value = 94524
$ws.Range("I11").Value = 3010
$ws.Range("M11").Formula = "=75000+13050"
$ws.Range("N11").Value = 2494

# ---------------------------------------------------------------
# Row 12  (15-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C12").Value = 0
$ws.Range("F12").Value = 60081
$ws.Range("I12").Value = 88
$ws.Range("M12").Formula = "=30000+29410"
$ws.Range("N12").Value = 583

# ---------------------------------------------------------------
# Row 13  (16-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C13").Value = 1738
$ws.Range("D13").Value = "SALAMI"
$ws.Range("F13").Value = 47537
$ws.Range("I13").Value = 15
$ws.Range("M13").Formula = "=20000+25440"
$ws.Range("N13").Value = 348

# ---------------------------------------------------------------
# Row 14  (17-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C14").Value = 4190
$ws.Range("D14").Value = "CHORIZO"
$ws.Range("F14").Value = 51862
$ws.Range("I14").Value = 40
$ws.Range("M14").Value = 45290
$ws.Range("N14").Value = 2345

# ---------------------------------------------------------------
# Row 15  (18-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C15").Value = 0
$ws.Range("F15").Value = 68774
$ws.Range("I15").Value = 15
$ws.Range("M15").Formula = "=500+43260+25000"

# ---------------------------------------------------------------
# Row 16  (19-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C16").Value = 752
$ws.Range("D16").Value = "TOCINO   "
$ws.Range("F16").Value = 71014
$ws.Range("I16").Value = 40
$ws.Range("M16").Formula = "=30000+39840"
$ws.Range("N16").Value = 385

# ---------------------------------------------------------------
# Row 17  (20-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = "ARABE"
$ws.Range("F17").Value = 70966
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = 44520
$ws.Range("K17").Value = "NOMINA # 47"
$ws.Range("L17").Value = 13341.34
$ws.Range("M17").Formula = "=30000+22220"
$ws.Range("N17").Value = 5300

# ---------------------------------------------------------------
# Row 18  (21-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C18").Value = 3983
$ws.Range("D18").Value = "RETAZO--TOCINO"
$ws.Range("F18").Value = 152652
$ws.Range("I18").Value = 0
$ws.Range("M18").Formula = "=70000+60000+17770"
$ws.Range("N18").Value = 903

# ---------------------------------------------------------------
# Row 19  (22-Nov-2021)
# ---------------------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = "CEBOLLA   "
$ws.Range("F19").Value = 47787
$ws.Range("I19").Value = 45
$ws.Range("M19").Formula = "=25000+22350"
$ws.Range("N19").Value = 383

# ---------------------------------------------------------------
# Recalculate so the Q column (variance = P - F) reflects the new
# entries before we touch its conditional font colour.
# ---------------------------------------------------------------
$excel.Calculate()

# The workbook has no real conditional formatting for column Q; the
# author manually painted the variance cells red when they go
# negative (bold, size 12, red) and black (bold, size 12) otherwise.
function Set-VarianceStyle($cell) {
    $v = $ws.Range($cell).Value2
    $ws.Range($cell).Font.Size = 12
    $ws.Range($cell).Font.Bold = $true
    if ($v -lt 0) {
        $ws.Range($cell).Font.Color = 255
    } else {
        $ws.Range($cell).Font.ThemeColor = 1
    }
}

Set-VarianceStyle "Q10"
Set-VarianceStyle "Q11"
Set-VarianceStyle "Q15"
Set-VarianceStyle "Q17"

# ---------------------------------------------------------------
# Leave the selection / scroll position where the author left it.
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("Q19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 7
